$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (rows 3 and 4),
# pushing the existing rows 3-9 down to rows 5-11.
$ws.Rows("3:4").Insert()

# Row 3: newest week's data (Primera)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44965
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107011
$ws.Range("J3").Value = "Tuna"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 34000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 34600
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 1922
$ws.Range("T3").Value = 18

# Row 4: newest week's data (Segunda)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44965
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 32000
$ws.Range("O4").Value = 33000
$ws.Range("P4").Value = 32333
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1796
$ws.Range("T4").Value = 18

# Ensure date cells keep the date number format used elsewhere in column D
$ws.Range("D3:D4").NumberFormat = $ws.Range("D5").NumberFormat
